$d = $word.ActiveDocument

# --- First paragraph: add a paragraph border, widen the left indent, and
#     fold the "ID token" run + trailing-space run into a single updated run.
$p1 = $d.Paragraphs(1)

# Adds <w:pBdr><w:top w:space="5"/><w:left w:space="5"/><w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

# <w:ind w:left="120"/> -> <w:ind w:left="225"/> (225 twips = 11.25 pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Replace the ID token together with the trailing space run that followed it,
# collapsing both runs into the single updated run the diff expects.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5309_topic_2__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5309__ID**", 2)
